# Actualización automática 2025-09-18 09:12:30
#
# Moves 43.1 of sales from "BRAVO MONTENEGRO DANIEL ANDRES" (row 3) to
# "FABIMP BENIGNO BRAVO S.A.S." (row 8), and records a new 122.67 / 13.81
# sale for "ILLER LOPEZ ROBERTO FERNANDO" (row 10) in the PORCELANATO /
# NO RESURTIBLES groups on the "VENTAS POR GRUPO" sheet, and the matching
# monthly figures on the "VENTA MENSUAL" sheet.

$wb = $excel.ActiveWorkbook

$wsGrupo  = $wb.Worksheets("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets("VENTA MENSUAL")

# --- VENTAS POR GRUPO (column M = PORCELANATO, column P = NO RESURTIBLES) ---
$wsGrupo.Range("M3").Value  = 0
$wsGrupo.Range("M8").Value  = 43.1
$wsGrupo.Range("M10").Value = 122.67
$wsGrupo.Range("P10").Value = 13.81

# Totals row ("X de 10" counters)
$wsGrupo.Range("M12").Value = "3 de 10"
$wsGrupo.Range("P12").Value = "1 de 10"

# --- VENTA MENSUAL (column F = septiembre) ---
$wsMensual.Range("F3").Value  = 0
$wsMensual.Range("F8").Value  = 43.1
$wsMensual.Range("F10").Value = 136.48

# Totals row (sum of F2:F11)
$wsMensual.Range("F12").Value = 5551.6
